{"js": "// The author renamed the character \"Zaq\" to \"Arcane\" throughout the story.\n// (Word's proofing engine also re-stamped nearby proofErr/run boundaries as a\n// side effect of the resave, but that carries no visible content change.)\nconst body = context.document.body;\n\nconst results = body.search(\"Zaq\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Arcane\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The author renamed the character \"Zaq\" to \"Arcane\" throughout the story.\n# (Word's proofing engine also re-stamped nearby proofErr/run boundaries as a\n# side effect of the resave, but that carries no visible content change.)\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Zaq\"\n$find.Replacement.Text = \"Arcane\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
